# Trade #81 closed at 2026-02-17 15:52:51 - unknown UNKNOWN +0.000%
#
# Updates the "Summary", "Strategy Status", "All Trades" and "MarketMaking"
# sheets to reflect the newly closed trade #81 (MarketMaking strategy).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Summary sheet
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1199.97   # Current Capital
$summary.Range("B4").Value = -0.04     # Total P&L $
$summary.Range("B6").Value = 81        # Total Trades
$summary.Range("B8").Value = 42        # Losing Trades
$summary.Range("B9").Value = 33.33     # Win Rate %

# ---------------------------------------------------------------------
# 2) Strategy Status sheet (MarketMaking row)
# ---------------------------------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C4").Value = 99.97      # Capital
$status.Range("D4").Value = 81         # Trades
$status.Range("E4").Value = -0.04      # P&L $
$status.Range("F4").Value = -0.03      # P&L %
$status.Range("G4").Value = 33.33      # Win Rate %

# ---------------------------------------------------------------------
# 3) Append the new trade row (#81) to both the "All Trades" and
#    "MarketMaking" sheets - they mirror each other.
# ---------------------------------------------------------------------
function Add-TradeRow($ws) {
    $row = 82

    $ws.Cells.Item($row, 1).Value = 81

    # Date-like text must be forced to text (leading apostrophe) so Excel
    # does not auto-convert it into a date serial number.
    $ws.Cells.Item($row, 2).Value = "'2026-02-17"
    $ws.Cells.Item($row, 3).Value = "15:52:45"

    $ws.Cells.Item($row, 4).Value = "MarketMaking"
    $ws.Cells.Item($row, 5).Value = "DOWN"
    $ws.Cells.Item($row, 6).Value = 0.21
    $ws.Cells.Item($row, 7).Value = 0.2
    $ws.Cells.Item($row, 8).Value = "CLOSED"
    $ws.Cells.Item($row, 9).Value = -4.7619
    $ws.Cells.Item($row, 10).Value = -0.01
    $ws.Cells.Item($row, 11).Value = 99.97
    $ws.Cells.Item($row, 12).Value = 0
    $ws.Cells.Item($row, 13).Value = 0
    $ws.Cells.Item($row, 14).Value = 0.6
    $ws.Cells.Item($row, 15).Value = "Normal spread capture: 19600 bps"
    $ws.Cells.Item($row, 16).Value = "early_exit"
    $ws.Cells.Item($row, 17).Value = 0.13
}

Add-TradeRow($wb.Worksheets.Item("All Trades"))
Add-TradeRow($wb.Worksheets.Item("MarketMaking"))
